$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make room for the new rows -------------------------------------------------
# Bad-driver block: one new row before the existing data rows, and two more
# new rows before the (soon to be shifted) Totals row.
$ws.Rows("3:3").Insert()
$ws.Rows("6:7").Insert()

# Good-driver block: one new row before the first data row (now at row 16
# after the three inserts above).
$ws.Rows("16:16").Insert()

# --- Bad Drivers table ----------------------------------------------------------
$ws.Range("A3").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 23.80.0.7"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 335
$ws.Range("D3").Value = 86

$ws.Range("A4").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.250.10.1"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 1893
$ws.Range("D4").Value = 93.2

$ws.Range("A5").Value = "Intel(R) Wi-Fi 6E AX211 160MHz - 22.160.0.3"
$ws.Range("B5").Value = 45
$ws.Range("C5").Value = 4832
$ws.Range("D5").Value = 98.40000000000001

$ws.Range("A6").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.100.0.3"
$ws.Range("B6").Value = 12
$ws.Range("C6").Value = 840
$ws.Range("D6").Value = 98.5

$ws.Range("A7").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.110.1.1"
$ws.Range("B7").Value = 64
$ws.Range("C7").Value = 6619
$ws.Range("D7").Value = 98.59999999999999

$ws.Range("B8").Value = 125
$ws.Range("C8").Value = 14519

# --- Good Drivers table ----------------------------------------------------------
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Range("B16").Value = 11128
$ws.Range("D16").Value = 100

# Row 16 inherited row 15's header formatting on insert; restore the "Total
# Samples" number style (bold, #,##0) used by the rest of the table by
# copying it down from the row right below.
$ws.Range("B17").Copy()
$ws.Range("B16").PasteSpecial(-4122)

$ws.Range("B17").Value = 486214

$ws.Range("B18").Value = 11140

$ws.Range("B19").Value = 14487

$ws.Range("B21").Value = 79953

$ws.Range("B22").Value = 35355

$ws.Range("B23").Value = 65425

$ws.Range("B24").Value = 117653
